$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = 0.0001201629638671875
$ws.Range("D2").Value = 0.4169491525423729
$ws.Range("F2").Value = 0.0001289844512939453

# Row 3
$ws.Range("B3").Value = 679
$ws.Range("C3").Value = 0.0001001358032226562
$ws.Range("D3").Value = 0.6630859375
$ws.Range("F3").Value = 0.0002388954162597656

# Row 4
$ws.Range("B4").Value = 33
$ws.Range("C4").Value = 0.00004982948303222656
$ws.Range("D4").Value = 0.9428571428571428
$ws.Range("F4").Value = 0.00003194808959960938

# Row 5
$ws.Range("C5").Value = 0.00005292892456054688
$ws.Range("F5").Value = 0.00005292892456054688

# Row 6
$ws.Range("B6").Value = 42
$ws.Range("C6").Value = 0.0001249313354492188
$ws.Range("D6").Value = 0.8076923076923077
$ws.Range("F6").Value = 0.00006580352783203125

# Row 7
$ws.Range("B7").Value = 79
$ws.Range("C7").Value = 0.00003814697265625
$ws.Range("D7").Value = 0.7383177570093458
$ws.Range("F7").Value = 0.0003659725189208984

# Row 8
$ws.Range("B8").Value = 9588
$ws.Range("C8").Value = 0.00007605552673339844
$ws.Range("D8").Value = 0.9816729804443535
$ws.Range("F8").Value = 0.0002660751342773438

# Row 9
$ws.Range("C9").Value = 0.00003314018249511719
$ws.Range("F9").Value = 0.00008511543273925781

# Row 10
$ws.Range("B10").Value = 923
$ws.Range("C10").Value = 0.0002419948577880859
$ws.Range("D10").Value = 0.9004878048780488
$ws.Range("F10").Value = 0.0006549358367919922

